$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A74").Value = "Golang Backend Developer-IoT Platforms-67189"
$ws.Range("B74").Value = "https://www.dice.com/job-detail/319af849-3254-43dc-9da8-7f0ff0fa9097"
$ws.Range("C74").Value = "Atlanta, Georgia"
$ws.Range("D74").Value = "Contract"
$ws.Range("E74").Value = '$$50/hr'
$ws.Range("F74").Value = "Robustware"

$ws.Range("A75").Value = "Golang Developer - Phoenix, AZ (Onsite)"
$ws.Range("B75").Value = "https://www.dice.com/job-detail/cde8083c-5502-4d06-8ff8-ad9d5cb597e2"
$ws.Range("C75").Value = "Phoenix, Arizona"
$ws.Range("D75").Value = "Contract, Third Party"
$ws.Range("E75").Value = "Depends on Experience"
$ws.Range("F75").Value = "STAND 8"
